# going-serverless-2017.pptx edit:
#   - slide 6:  title "Patterns" -> "Use Cases" (2nd title line unchanged)
#               last content bullet "And more..." -> two new bullets:
#                 "Mobile backends" and "IoT backends"
#   - slide 14: title "Patterns" / "again" -> "Use Case" / "patterns"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 6: "Patterns / we will see this again" -> "Use Cases / ..."
# ---------------------------------------------------------------
$slide6 = $p.Slides.Item(6)

$title6 = $slide6.Shapes.Title.TextFrame.TextRange
# Title text box has two runs: "Patterns" then (after a line break) "we
# will see this again". Only the first run's text changes.
$title6.Runs(1, 1).Text = "Use Cases"

# Content placeholder: replace the final "And more..." bullet with two
# new bullets describing concrete use-case patterns.
$body6 = $slide6.Shapes.Placeholders.Item(2).TextFrame.TextRange
$lastBullet6 = $body6.Paragraphs(5)
$lastBullet6.Text = "Mobile backends"
$lastBullet6.InsertAfter("`rIoT backends") | Out-Null

# ---------------------------------------------------------------
# Slide 14: "Patterns / again" -> "Use Case / patterns"
# ---------------------------------------------------------------
$slide14 = $p.Slides.Item(14)

$title14 = $slide14.Shapes.Title.TextFrame.TextRange
$title14.Runs(1, 1).Text = "Use Case"
$title14.Runs(2, 1).Text = "patterns"
